$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Based on results from Q10 and Q15, ..." -> "... Q11 and Q15, ..."
#    (Q16's explanatory paragraph references Q10 by mistake; fix to Q11.)
#    Scope the Find/Replace tightly to that one paragraph so the unrelated
#    "Q10." heading (different question) is left untouched.
# ---------------------------------------------------------------------------
$q16Body = $d.Content
$q16Body.Find.Execute("Based on results from", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$q16ParaRange = $q16Body.Paragraphs(1).Range
$scoped = $d.Range($q16ParaRange.Start, $q16ParaRange.End)
$scoped.Find.Execute("Q10", $true, $true, $false, $false, $false, $true, 0, $false, "Q11", 2)

# ---------------------------------------------------------------------------
# 2) Insert a new "Challenge Problems" (Heading 2) paragraph right before the
#    existing "Q17. (CHALLENGE PROBLEM)" heading, and wrap it (through the
#    end of the document) in a new "challenge-problems" bookmark.
# ---------------------------------------------------------------------------
$q16Body2 = $d.Content
$q16Body2.Find.Execute("Based on results from", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$q16ParaRange2 = $q16Body2.Paragraphs(1).Range
$insertAt = $q16ParaRange2.End
$q16ParaRange2.InsertParagraphAfter()

$newHeadingRange = $d.Range($insertAt, $insertAt)
$newHeadingPara = $newHeadingRange.Paragraphs(1)
$newHeadingPara.Range.Text = "Challenge Problems"
$newHeadingPara.Style = $d.Styles.Item("Heading 2")

$bmRange = $d.Range($insertAt, $d.Content.End)
$d.Bookmarks.Add("challenge-problems", $bmRange)

Write-Output "done"
